# Apply updated symbol list values (price/volume/coin refreshes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '254.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '3.64%'
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '28.04'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-3.93%'
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.355'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '4.08%'
$ws.Range('E4').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.54%'
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.712'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.03%'
$ws.Range('E6').Style = 'Normal'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('B7').Style = 'Normal'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C7').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8630'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '0.73%'
$ws.Range('E7').Style = 'Normal'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('B8').Style = 'Normal'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('C8').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9123'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '5.69%'
$ws.Range('E8').Style = 'Normal'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('B9').Style = 'Normal'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('C9').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1423'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '3.59%'
$ws.Range('E9').Style = 'Normal'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('B10').Style = 'Normal'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('C10').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07191'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '1.66%'
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('B11').Style = 'Normal'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('C11').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.03206'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-2.45%'
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09240'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-1.37%'
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001540'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.33%'
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'One'
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0006060'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.94%'
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.005885'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-2.28%'
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'LEO'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.496'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.26%'
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.226'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.60%'
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.252'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '3.22%'
$ws.Range('E18').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '2.91%'
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1304'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '1.74%'
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.528'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '10.93%'
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04152'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '0.20%'
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-1.39%'
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.005120'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '23.59%'
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.001225'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-0.01%'
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-0.72%'
$ws.Range('E27').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03841'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '3.02%'
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1100'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2.74%'
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'KickToken'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.002950'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-49.17%'
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002380'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '8.30%'
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.009998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '9.07%'
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005290'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.05%'
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.12%'
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1000'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '72.60%'
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002213'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '1.87%'
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.12%'
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.12%'
$ws.Range('E50').Style = 'Normal'

Write-Output "Applied 91 cell updates"
